# Add excel export to linear model
# Adds a new "ini_2" row (row 7) below the existing "ini" row (row 6),
# with integer-formatted values B7:K7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inst1")

# Label for the new row
$ws.Range("A7").Value = "ini_2"

# New data values (integers) for B7:K7
$values = @(6, 10, 6, 7, 8, 8, 9, 12, 9, 10)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i  # B = 2
    $cell = $ws.Cells.Item(7, $col)
    $cell.Value = $values[$i]
    $cell.NumberFormat = "0"
}

# Update selection to match the recorded end-state
$ws.Range("I11").Select()
